$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# New localization-report row data (the handoff that "failed transform")
# ---------------------------------------------------------------------------
$oldFile    = "f6d729bf-dd8f-4bdc-ab89-d264195cab48.md"
$newFile    = "154c1f0d-5526-4ddb-8004-e9dd61f70998.md"
$failFile   = "76f1699a-91da-4944-b87d-d96f2df952ae.md"
$failStatus = "Handoff transform failed"
$epoch      = "0001-01-01 00:00:00"
$ignored    = "Ignored"

$zhOldXlf = "f6d729bf-dd8f-4bdc-ab89-d264195cab48.dd7e42755ba01a6ec27b46a0f786d290dba0a38a.zh-cn.xlf"
$zhNewXlf = "154c1f0d-5526-4ddb-8004-e9dd61f70998.5d2089412d93602adcd264b3f59baf1f4cb7f9f9.zh-cn.xlf"
$deOldXlf = "f6d729bf-dd8f-4bdc-ab89-d264195cab48.dd7e42755ba01a6ec27b46a0f786d290dba0a38a.de-de.xlf"
$deNewXlf = "154c1f0d-5526-4ddb-8004-e9dd61f70998.5d2089412d93602adcd264b3f59baf1f4cb7f9f9.de-de.xlf"

$zhOldDt = "2016-01-11 05:23:26"
$zhNewDt = "2016-01-11 05:24:32"
$deOldDt = "2016-01-11 05:23:42"
$deNewDt = "2016-01-11 05:24:48"

$baseUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/14876332aac6ca5a5dc3fe181950b268581e5fbf"
$newMdUrl  = "$baseUrl/e2e/$newFile"
$failMdUrl = "$baseUrl/e2e/$failFile"
$cfgUrl    = "$baseUrl/.localization-config"

# ===========================================================================
# Sheet "Overview"
# ===========================================================================
$ws = $wb.Worksheets.Item("Overview")

# Make room for the new row between the "Ready for handoff" row and the
# ".localization-config" row.
$ws.Rows.Item(3).Insert()

$ws.Range("A2").Value = $newFile

$ws.Range("A3").Value = $failFile
$ws.Range("B3").Value = $failStatus
$ws.Range("C3").Value = $failStatus

# Row 4 keeps the original ".localization-config" / "Not to be localized"
# content (it just shifted down one row) - nothing to change there.

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $newMdUrl, "", "", $newFile) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), $failMdUrl, "", "", $failFile) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), $cfgUrl, "", "", ".localization-config") | Out-Null

# ===========================================================================
# Sheet "zh-cn"
# ===========================================================================
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Rows.Item(3).Insert()

$ws.Range("A2").Value = $newFile
$ws.Range("C2").Value = $zhNewXlf
$ws.Range("D2").Value = $zhNewDt

$ws.Range("A3").Value = $failFile
$ws.Range("B3").Value = $failStatus
$ws.Range("C3").Clear()
$ws.Range("D3").Value = $epoch
$ws.Range("G3").Value = $epoch
$ws.Range("H3").Value = $ignored

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $newMdUrl, "", "", $newFile) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/236bd29ca5f608564d792c26dd4bd69a29308513/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/$zhNewXlf", "", "", $zhNewXlf) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), $failMdUrl, "", "", $failFile) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), $cfgUrl, "", "", ".localization-config") | Out-Null

# ===========================================================================
# Sheet "de-de"
# ===========================================================================
$ws = $wb.Worksheets.Item("de-de")

$ws.Rows.Item(3).Insert()

$ws.Range("A2").Value = $newFile
$ws.Range("C2").Value = $deNewXlf
$ws.Range("D2").Value = $deNewDt

$ws.Range("A3").Value = $failFile
$ws.Range("B3").Value = $failStatus
$ws.Range("C3").Clear()
$ws.Range("D3").Value = $epoch
$ws.Range("G3").Value = $epoch
$ws.Range("H3").Value = $ignored

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $newMdUrl, "", "", $newFile) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c0296cdd6c7eda395f9458585b0d3a812e3ee5da/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/$deNewXlf", "", "", $deNewXlf) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), $failMdUrl, "", "", $failFile) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), $cfgUrl, "", "", ".localization-config") | Out-Null
